$wb = $excel.ActiveWorkbook

# 1. Rename header on "Weekly Quantity" sheet
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header on "Monthly Trend" sheet
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add new "PO Forecast" sheet after "Monthly Trend"
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Re-fetch Monthly Trend reference (stale after Add) then move Forecast after it
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $wsMonthly)

# Re-fetch Forecast reference again since Move() invalidates previous handle
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other sheets: bold, thin border,
# centered/top-aligned text.
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$data = @(
    @(45144.99999999999, 0, -74.11465617871889, 48.67959643664039),
    @(45158.99999999999, 0, -70.58070690887989, 54.87595740809189),
    @(45179.99999999999, 0, -66.6214458109021, 59.5608444602473),
    @(45207.99999999999, 0, -55.82291843484577, 63.56590573070697),
    @(45228.99999999999, 4, -58.37911715636312, 66.37387467691147),
    @(45277.99999999999, 13, -44.59190065701024, 73.89749386006662),
    @(45305.99999999999, 18, -43.75828306126981, 78.23459939438651),
    @(45319.99999999999, 21, -38.94472521028322, 82.80456205912095),
    @(45326.99999999999, 22, -41.56529147944974, 80.22992009847869),
    @(45333.99999999999, 23, -40.16069539635289, 81.94672175092307),
    @(45354.99999999999, 27, -33.14915674219721, 88.70209759137802),
    @(45361.99999999999, 29, -31.16921571386285, 89.12422081284257),
    @(45375.99999999999, 31, -29.51984846890462, 87.58680835350704),
    @(45396.99999999999, 35, -31.2690173884736, 94.52106229521092),
    @(45410.99999999999, 38, -23.34196240160038, 99.66628297781405),
    @(45417.99999999999, 39, -22.2855178435832, 98.66510609714426),
    @(45424.99999999999, 41, -19.95184111293637, 100.8285168187325),
    @(45431.99999999999, 42, -20.31056104194815, 97.79606343474815),
    @(45438.99999999999, 43, -16.46572504193887, 102.4024246892847),
    @(45445.99999999999, 44, -12.11098058286826, 104.3851169832152),
    @(45452.99999999999, 46, -13.02477139592461, 105.0500899367),
    @(45459.99999999999, 47, -14.54958821329656, 103.789938087363),
    @(45487.99999999999, 52, -10.17643517398805, 110.9801004658706),
    @(45494.99999999999, 54, -7.166030731771311, 109.619312577664),
    @(45501.99999999999, 55, -5.258766461379471, 113.1731786341769),
    @(45508.99999999999, 56, -0.6087987581460336, 117.2898459429282),
    @(45515.99999999999, 58, -0.4623303532139852, 119.5142885064067),
    @(45522.99999999999, 59, 0.5615523902777458, 117.3777887507175),
    @(45529.99999999999, 60, 1.291341039404485, 121.7768219789893),
    @(45536.99999999999, 62, 2.694356681050304, 123.0514796926838),
    @(45543.99999999999, 63, 2.280538436849357, 120.1385197000369)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Column A holds dates; give it the same date/time number format used for
# the "ds" columns on the other two sheets.
$wsForecast.Range("A2:A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Match the page margins used on the other two sheets.
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36
